# Add a "Workflow State" column (column I) to the exemplar link-checker
# report sheet: header in I1, "man!" filler values in I2:I10, matching
# the formatting of the existing last column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell
$ws.Range("I1").Value = "Workflow State"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-10
For ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 9).Value = "man!"
}
$ws.Range("H2").Copy()
$ws.Range("I2:I10").PasteSpecial(-4122)  # xlPasteFormats

[void]($excel.CutCopyMode = $false)

# Match the recorded selection from the diff
[void]$ws.Range("H8").Select()
